$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 11
$ws.Range("G5").Value = 5
